# Updated cryptos list on Mon Apr 24 12:09:31 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto-ranking table, and fixes the ordering of two coin pairs whose rows
# had their Coin/Link/Price/Volume values swapped (TrustWalletToken <->
# FraxShare at rows 39/40, and Decentraland <-> WEMIXTOKEN at rows 45/46).
#
# Many of the new Price strings (e.g. "1.010", "45.39") are valid-looking
# numbers, but the sheet stores them as plain text (no leading "="), so a
# naive `.Value = "1.010"` assignment would get auto-coerced to the number
# 1.01 and lose the trailing zero. To keep them as literal text without
# permanently tattooing the cell with a new NumberFormat/style, we briefly
# enter them with a leading apostrophe (forces text entry, same as typing
# '1.010 into Excel) and immediately reset the cell back to the builtin
# "Normal" style so no stray style index is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.752.57'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.870.57'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("D4").Value = '''1.010'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").Value = '''335.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '
$ws.Range("D6").Value = '''1.010'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").Value = '''0.4697'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").Value = '''0.08011'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = '''45.39'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.37%  '
$ws.Range("D11").Value = '''1.006'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").Value = '''21.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.66%  '
$ws.Range("D13").Value = '1.872.02'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").Value = '''6.016'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '''7.267'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").Value = '''1.012'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = '''88.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("D18").Value = '''0.06748'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.58%  '
$ws.Range("D19").Value = '''0.00001045'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("D20").Value = '''17.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("D21").Value = '''1.010'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").Value = '27.747.54'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").Value = '''5.490'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("D24").Value = '''10.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("D25").Value = '''2.324'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").Value = '2.092.96'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '''159.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D29").Value = '''2.161'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.40%  '
$ws.Range("D30").Value = '''5.463'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("D31").Value = '''121.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").Value = '''0.9809'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("D33").Value = '''0.09489'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").Value = '''3.615'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").Value = '''5.327'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = '''1.343'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.57%  '
$ws.Range("D37").Value = '''0.06072'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.02%  '
$ws.Range("D38").Value = '''0.02239'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '''1.200'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.51%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''8.321'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.19%  '
$ws.Range("E41").Value = '  +0.40%  '
$ws.Range("D42").Value = '''0.5995'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.31%  '
$ws.Range("D43").Value = '''0.1891'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").Value = '''10.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.5669'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '''1.241'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = '''0.06763'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.17%  '
$ws.Range("D50").Value = '''112.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.61%  '
$ws.Range("D51").Value = '''3.045'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -10.11%  '
